# Implement viral load indicators
# Adds a new indicator block (row 37 header cells) mirroring the existing
# "HEPATITIS_B_DIAGNOSIS_DURING_PERIOD" header row (row 36), re-purposed for
# the new "HAD_CV_WITHIN_12_MONTH_AFTER_ARV_START" indicator, and adjusts
# row heights / selection state that shifted as a result of the edit.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("File active")
$ws2 = $wb.Worksheets.Item("Schémas ARV")

# --- Row 37: fill in the new indicator's sub-column headers ---------------
$ws1.Range("B37").Value = "{key: 'HAD_CV_WITHIN_12_MONTH_AFTER_ARV_START', gender: 0}"
$ws1.Range("C37").Value = "{key: 'HAD_CV_WITHIN_12_MONTH_AFTER_ARV_START', gender: 1}"
$ws1.Range("D37").Value = "{key:'HAD_CV_WITHIN_12_MONTH_AFTER_ARV_START', age_max: 15}"
$ws1.Range("E37").Value = "{key: 'HAD_CV_WITHIN_12_MONTH_AFTER_ARV_START', age_min: 15}"
$ws1.Range("F37").Value = "{key: 'HAD_CV_WITHIN_12_MONTH_AFTER_ARV_START'}"
$ws1.Range("G37").Value = "{key: 'HAD_CV_WITHIN_12_MONTH_AFTER_ARV_START', age_max: 15, gender: 0}"
$ws1.Range("H37").Value = "{key: 'HAD_CV_WITHIN_12_MONTH_AFTER_ARV_START', age_max: 15, gender: 1}"
$ws1.Range("I37").Value = "{key: 'HAD_CV_WITHIN_12_MONTH_AFTER_ARV_START', age_min: 15, gender: 0}"
$ws1.Range("J37").Value = "{key: 'HAD_CV_WITHIN_12_MONTH_AFTER_ARV_START', age_min: 15, gender: 1}"
$ws1.Range("K37").Value = "{key: 'HAD_CV_WITHIN_12_MONTH_AFTER_ARV_START', age_max: 1, gender: 0}"
$ws1.Range("L37").Value = "{key: 'HAD_CV_WITHIN_12_MONTH_AFTER_ARV_START', age_max: 1, gender: 1}"
$ws1.Range("M37").Value = "{key: 'HAD_CV_WITHIN_12_MONTH_AFTER_ARV_START', age_min: 1, age_max: 4, gender: 0}"
$ws1.Range("N37").Value = "{key: 'HAD_CV_WITHIN_12_MONTH_AFTER_ARV_START', age_max: 1, gender: 0}"
$ws1.Range("O37").Value = "{key: 'HAD_CV_WITHIN_12_MONTH_AFTER_ARV_START', age_min: 5, age_max: 9, gender: 0}"
$ws1.Range("P37").Value = "{key: 'HAD_CV_WITHIN_12_MONTH_AFTER_ARV_START', age_min: 5, age_max: 9, gender: 1}"
$ws1.Range("Q37").Value = "{key: 'HAD_CV_WITHIN_12_MONTH_AFTER_ARV_START', age_min: 10, age_max: 14, gender: 0}"
$ws1.Range("R37").Value = "{key: 'HAD_CV_WITHIN_12_MONTH_AFTER_ARV_START', age_min: 10, age_max: 14, gender: 1}"
$ws1.Range("S37").Value = "{key: 'HAD_CV_WITHIN_12_MONTH_AFTER_ARV_START', age_min: 15, age_max: 19, gender: 0}"
$ws1.Range("T37").Value = "{key: 'HAD_CV_WITHIN_12_MONTH_AFTER_ARV_START', age_min: 15, age_max: 19, gender: 1}"
$ws1.Range("U37").Value = "{key: 'HAD_CV_WITHIN_12_MONTH_AFTER_ARV_START', age_min: 20, age_max: 24, gender: 0}"
$ws1.Range("V37").Value = "{key: 'HAD_CV_WITHIN_12_MONTH_AFTER_ARV_START', age_min: 20, age_max: 24, gender: 1}"
$ws1.Range("W37").Value = "{key: 'HAD_CV_WITHIN_12_MONTH_AFTER_ARV_START', age_min: 25, age_max: 49, gender: 0}"
$ws1.Range("X37").Value = "{key: 'HAD_CV_WITHIN_12_MONTH_AFTER_ARV_START', age_min: 25, age_max: 49, gender: 1}"
$ws1.Range("Y37").Value = "{key: 'HAD_CV_WITHIN_12_MONTH_AFTER_ARV_START', age_min: 50, gender: 0}"
$ws1.Range("Z37").Value = "{key: 'HAD_CV_WITHIN_12_MONTH_AFTER_ARV_START', age_min: 50, gender: 1}"
$ws1.Range("AA37").Value = "{key: 'HAD_CV_WITHIN_12_MONTH_AFTER_ARV_START', gender: 0}"
$ws1.Range("AB37").Value = "{key: 'HAD_CV_WITHIN_12_MONTH_AFTER_ARV_START', gender: 1}"

# Row 37 now has headers like row 36 above it (which carries the matching
# formatting for this kind of sub-header row) - copy that formatting over so
# the new cells get the same font/fill/border/alignment/number format.
$ws1.Range("B36:AB36").Copy()
$ws1.Range("B37:AB37").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row heights that changed to fit the new/adjusted content -------------
$ws1.Rows.Item(37).RowHeight = 53.7
$ws1.Rows.Item(38).RowHeight = 40.25
$ws1.Rows.Item(39).RowHeight = 56.7

# --- Selection / cursor state ----------------------------------------------
# (Select on sheet2 first, then sheet1 last so "File active" ends up as the
# active/visible tab again, matching the saved view state.)
$ws2.Range("A1").Select()
$ws1.Range("B37:AB37").Select()
